$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two observation rows (2 and 3) got their record identity swapped:
# row 2's Id/Antal/Ost/Nord/Starttid/Sluttid moved to row 3 and vice versa.
# All the other columns already hold identical values between the two rows,
# so only these six columns need to be exchanged.
#
# A scratch cell (AZ1) outside the used range (A1:AY3) is used to hold a
# value temporarily while swapping a pair of cells via copy/paste, which
# preserves the original cell type/formatting (e.g. text "100" in column I
# instead of it being re-interpreted as a number).

function Swap-CellValues {
    param($Sheet, [string]$RefTop, [string]$RefBottom, [string]$ScratchRef)

    $scratch = $Sheet.Range($ScratchRef)
    $top = $Sheet.Range($RefTop)
    $bottom = $Sheet.Range($RefBottom)

    # top -> scratch
    $top.Copy()
    $scratch.Value2 = $null
    $scratch.PasteSpecial(-4104)

    # bottom -> top
    $bottom.Copy()
    $top.Value2 = $null
    $top.PasteSpecial(-4104)

    # scratch -> bottom
    $scratch.Copy()
    $bottom.Value2 = $null
    $bottom.PasteSpecial(-4104)

    $scratch.Clear()
}

$scratchCell = "AZ1"

Swap-CellValues $ws "A2"  "A3"  $scratchCell   # Id
Swap-CellValues $ws "I2"  "I3"  $scratchCell   # Antal
Swap-CellValues $ws "Q2"  "Q3"  $scratchCell   # Ost
Swap-CellValues $ws "R2"  "R3"  $scratchCell   # Nord
Swap-CellValues $ws "Z2"  "Z3"  $scratchCell   # Starttid
Swap-CellValues $ws "AB2" "AB3" $scratchCell   # Sluttid
